$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $ok = $r.Find.Execute($oldText, $false, $false, $false, $false, $false, `
                           $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Replace failed in paragraph $paraIndex for text: $oldText"
    }
}

# 1) Title: "School of Mines" -> "Computer Science" (keep the trailing " Worksheet")
Replace-InParagraph 1 "School of Mines" "Computer Science"

# 2) Question list items (each is its own single-paragraph ListParagraph item)
Replace-InParagraph 3  "What is the School of Mines' mission?" "What is computer science? "
Replace-InParagraph 7  "How many undergraduate and graduate degree programs are offered at the School of Mines?" "What are the main goals of computer science? "
Replace-InParagraph 11 "What are the School of Mines' research priorities?" "What is the difference between computer science and information technology? "
Replace-InParagraph 15 "What are the School of Mines' teaching priorities?" "What are the different types of computer systems? "
Replace-InParagraph 19 "How many faculty members are at the School of Mines?" "What is a programming language? "
Replace-InParagraph 23 "What are the School of Mines' facilities and resources?" "How do computers store and process information? "
Replace-InParagraph 27 "What is the admissions process for students seeking admission to the School of Mines?" "How do computers make decisions? "
Replace-InParagraph 31 "What type of student is best suited for study at the School of Mines?" "What are some common applications of computer science?"

# 3) Remove the last two question blocks entirely ("How much does it cost..." and
#    "What kind of career opportunities...") along with their tables and the blank
#    paragraph that separates them, while keeping one trailing blank paragraph.
#
# Layout right now (after the edits above), counting from the paragraph that now
# reads "What are some common applications of computer science?":
#   P(31) "What are some common applications of computer science?"
#   Table (immediately after P31)
#   P(32) blank
#   P(33) "How much does it cost to attend the School of Mines?"
#   Table
#   P(34) blank
#   P(35) "What kind of career opportunities are available to graduates of the School of Mines?"
#   Table
#   P(36) blank  <-- keep this one (final trailing paragraph)

$startPara = $d.Paragraphs(32)
$startPos = $startPara.Range.Start

$lastTable = $d.Tables($d.Tables.Count)
$endPos = $lastTable.Range.End

$killRange = $d.Range($startPos, $endPos)
$killRange.Delete()
